$d = $word.ActiveDocument

# The Pearson / BTEC logo pictures that live in this document's headers
# and footers were saved with mismatched "name" metadata - the two
# Pearson logo pictures (in the footers) were both called "image1.png"
# and the BTec logo picture (in the header) was called "image2.jpg".
# Rename each inline picture back to its correct, distinct file-based
# name:
#   Pearson logo pictures (footers)  image1.png -> image2.png
#   BTec logo picture     (header)   image2.jpg -> image1.jpg
#
# NOTE: InlineShapes whose Range spans more paragraphs than just their
# own (e.g. a logo that isn't the first paragraph of its header/footer)
# can end up with a stale anchor once they are fetched from the
# section's full Header/Footer Range, so re-fetch each shape through its
# own (single-paragraph) Range before renaming it - that keeps the
# rename targeted at exactly the right picture.

function Rename-LogoPicture($range, [string]$oldAltText, [string]$newName) {
    $shapes = $range.InlineShapes
    for ($k = 1; $k -le $shapes.Count; $k++) {
        $shp = $shapes.Item($k)
        if ($shp.AlternativeText -eq $oldAltText) {
            $self = $shp.Range.InlineShapes.Item(1)
            $self.Name = $newName
        }
    }
}

for ($s = 1; $s -le $d.Sections.Count; $s++) {
    $sec = $d.Sections.Item($s)

    for ($i = 1; $i -le 3; $i++) {
        $hdr = $sec.Headers.Item($i)
        if ($hdr.Exists) {
            Rename-LogoPicture $hdr.Range "BTec_Logo-Orange" "image1.jpg"
        }

        $ftr = $sec.Footers.Item($i)
        if ($ftr.Exists) {
            Rename-LogoPicture $ftr.Range "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" "image2.png"
        }
    }
}
